# Cập nhật lại phân công công việc, ghép công việc Lập Tuấn
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Update person in charge for row 8 (task "Xây dựng module quản lý dịch vụ, hóa đơn")
# from "Lập" to "Lập, Tuấn" (ghép công việc Lập + Tuấn)
$ws.Range("C8").Value = "Lập, Tuấn"

# Update the view state: top-left visible cell and active selection
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select()
